# Daily attendance processing - rotate the "Recorded By" (column G) list so
# that the first contributor listed moves to the end of the comma-separated
# sequence (e.g. "a, b, c" -> "b, c, a"). Cells with a single value, or whose
# value is exactly "admin@admin.com, System", are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells(1,1).SpecialCells(11).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -eq "admin@admin.com, System") {
        continue
    }

    $rawParts = $val.Split(",")
    if ($rawParts.Length -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $first = $parts[0]
    $rest = $parts[1..($parts.Length - 1)]
    $newParts = $rest + $first
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
